# GammaFiber2F-HW15 "AveragedIntensites" sheet: re-ran the averaging code to add
# the Gaussian-Quadrature scheme alongside 3 new spiral sampling schemes
# (Spiral-90deg-10rot-5space, Spiral-90deg-15rot-5space, Spiral-90deg-10rot-3space).
# This pushes the previously appended rows (NoRotation-tilt60deg .. HexGrid-60degTilt5degRes)
# down so the sheet now spans A1:M19 instead of A1:M16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column A's bold/bordered style (same as existing rows) down to the new rows.
$ws.Range("A14:A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.36567837976243
$ws.Range("D10").Value = 0.3828158769166418
$ws.Range("E10").Value = 1.050977968836757
$ws.Range("F10").Value = 1.36567837976243
$ws.Range("G10").Value = 0.6858759666533483
$ws.Range("H10").Value = 1.133475496403593
$ws.Range("I10").Value = 1.13367014862962
$ws.Range("J10").Value = 0.3828158769166418
$ws.Range("K10").Value = 0.7168969228766995
$ws.Range("L10").Value = 1.041287651319565
$ws.Range("M10").Value = 0.958748972867065

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9765742537872838
$ws.Range("D11").Value = 0.8508376388969759
$ws.Range("E11").Value = 1.070428173496901
$ws.Range("F11").Value = 0.9765742537872838
$ws.Range("G11").Value = 0.8692412970665807
$ws.Range("H11").Value = 1.263320494963714
$ws.Range("I11").Value = 1.029135551878591
$ws.Range("J11").Value = 0.8508376388969759
$ws.Range("K11").Value = 0.9606329061969385
$ws.Range("L11").Value = 0.9686035799921112
$ws.Range("M11").Value = 1.009922901681674

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9746220448296159
$ws.Range("D12").Value = 0.8522888801302883
$ws.Range("E12").Value = 1.070690758119042
$ws.Range("F12").Value = 0.9746220448296159
$ws.Range("G12").Value = 0.8702029211470393
$ws.Range("H12").Value = 1.262646001209882
$ws.Range("I12").Value = 1.028662307964872
$ws.Range("J12").Value = 0.8522888801302883
$ws.Range("K12").Value = 0.9614898191246652
$ws.Range("L12").Value = 0.9680559319771405
$ws.Range("M12").Value = 1.009852152233457

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9767318454647269
$ws.Range("D13").Value = 0.8505263712862053
$ws.Range("E13").Value = 1.070810500218904
$ws.Range("F13").Value = 0.9767318454647269
$ws.Range("G13").Value = 0.8692184703559986
$ws.Range("H13").Value = 1.261613160905498
$ws.Range("I13").Value = 1.029094525074254
$ws.Range("J13").Value = 0.8505263712862053
$ws.Range("K13").Value = 0.9606684357525548
$ws.Range("L13").Value = 0.9687001406086408
$ws.Range("M13").Value = 1.009665812217598

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.5465840000000001
$ws.Range("D14").Value = 0.4161720000000006
$ws.Range("E14").Value = 1.322248000000002
$ws.Range("F14").Value = 0.5465840000000001
$ws.Range("G14").Value = 0.4088119999999999
$ws.Range("H14").Value = 3.306859999999999
$ws.Range("I14").Value = 1.160432000000001
$ws.Range("J14").Value = 0.4161720000000006
$ws.Range("K14").Value = 0.8692100000000011
$ws.Range("L14").Value = 0.7078970000000007
$ws.Range("M14").Value = 1.193518

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.2
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 1.54
$ws.Range("F15").Value = 0.2
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 5.08555
$ws.Range("I15").Value = 1.28
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0.77
$ws.Range("L15").Value = 0.485
$ws.Range("M15").Value = 1.350925

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.5456073000960006
$ws.Range("D16").Value = 0.4141347309568011
$ws.Range("E16").Value = 1.322606180147202
$ws.Range("F16").Value = 0.5456073000960006
$ws.Range("G16").Value = 0.4101720221696006
$ws.Range("H16").Value = 3.307962592767999
$ws.Range("I16").Value = 1.160864543129604
$ws.Range("J16").Value = 0.4141347309568011
$ws.Range("K16").Value = 0.8683704555520013
$ws.Range("L16").Value = 0.7069888778240009
$ws.Range("M16").Value = 1.193557894877868

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9880549741619945
$ws.Range("D17").Value = 0.992829067359612
$ws.Range("E17").Value = 0.9953528709887164
$ws.Range("F17").Value = 0.9880549741619945
$ws.Range("G17").Value = 0.9887546906017016
$ws.Range("H17").Value = 0.9972423266463194
$ws.Range("I17").Value = 0.9925918642535275
$ws.Range("J17").Value = 0.992829067359612
$ws.Range("K17").Value = 0.9940909691741642
$ws.Range("L17").Value = 0.9910729716680794
$ws.Range("M17").Value = 0.9924709656686453

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9466312067446748
$ws.Range("D18").Value = 1.078602224028307
$ws.Range("E18").Value = 0.9959429379799353
$ws.Range("F18").Value = 0.9466312067446748
$ws.Range("G18").Value = 1.028897430557631
$ws.Range("H18").Value = 1.005351800859394
$ws.Range("I18").Value = 0.971319758758452
$ws.Range("J18").Value = 1.078602224028307
$ws.Range("K18").Value = 1.037272581004121
$ws.Range("L18").Value = 0.9919518938743979
$ws.Range("M18").Value = 1.004457559821399

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9783807829387059
$ws.Range("D19").Value = 1.250451318232327
$ws.Range("E19").Value = 0.926489151778749
$ws.Range("F19").Value = 0.9783807829387059
$ws.Range("G19").Value = 1.114291107181111
$ws.Range("H19").Value = 0.7886986867608025
$ws.Range("I19").Value = 0.9311999269302119
$ws.Range("J19").Value = 1.250451318232327
$ws.Range("K19").Value = 1.088470235005538
$ws.Range("L19").Value = 1.033425508972122
$ws.Range("M19").Value = 0.9982518289703178
